$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Rewrite the results table (new benchmark run + reordered rows)
#    Column A = trace name, B..F = predictor misprediction rates
# -----------------------------------------------------------------
$ws.Range("A2").Value = "mm_2"
$ws.Range("B2").Value = 8.4830000000000005
$ws.Range("C2").Value = 10.138
$ws.Range("D2").Value = 6.7519999999999998
$ws.Range("E2").Value = 6.6539999999999999
$ws.Range("F2").Value = 6.9349999999999996

$ws.Range("A3").Value = "mm_1"
$ws.Range("B3").Value = 2.581
$ws.Range("C3").Value = 6.6959999999999997
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1.0609999999999999
$ws.Range("F3").Value = 1.0029999999999999

$ws.Range("A4").Value = "int_2"
$ws.Range("B4").Value = 0.42599999999999999
$ws.Range("C4").Value = 0.42
$ws.Range("D4").Value = 0.27600000000000002
$ws.Range("E4").Value = 0.27900000000000003
$ws.Range("F4").Value = 0.27500000000000002

$ws.Range("A5").Value = "int_1"
$ws.Range("B5").Value = 12.622
$ws.Range("C5").Value = 13.839
$ws.Range("D5").Value = 10.753
$ws.Range("E5").Value = 10.77
$ws.Range("F5").Value = 11.407999999999999

$ws.Range("A6").Value = "fp_2"
$ws.Range("B6").Value = 3.246
$ws.Range("C6").Value = 1.6779999999999999
$ws.Range("D6").Value = 1.327
$ws.Range("E6").Value = 0.219
$ws.Range("F6").Value = 1.327

$ws.Range("A7").Value = "fp_1"
$ws.Range("B7").Value = 0.99099999999999999
$ws.Range("C7").Value = 0.82499999999999996
$ws.Range("D7").Value = 0.81
$ws.Range("E7").Value = 0.81100000000000005
$ws.Range("F7").Value = 0.81200000000000006

# -----------------------------------------------------------------
# 2. Worksheet cosmetics: new column width, selection, page setup
# -----------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 95

$ws.Range("E6").Select()

# -----------------------------------------------------------------
# 3. Chart: update title, add axis titles, tweak tick marks, resize
# -----------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$chart.ChartTitle.Text = "Performance of different branch predictors "

$catAx = $chart.Axes(1)
$catAx.HasTitle = $true
$catAx.AxisTitle.Text = "Trace"
$catAx.MajorTickMark = 3

$valAx = $chart.Axes(2)
$valAx.HasTitle = $true
$valAx.AxisTitle.Text = "Misprediction Rate (%)"
$valAx.MajorTickMark = 3

# Resize / reposition the chart to its new anchor
$co.Left = 335.4999
$co.Top = 209.5
$co.Width = 709.7501
$co.Height = 522.5
